$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Table 1 (Statistical Institution): Enterprises density (per 1000 people)
$ws.Range("B11").Value = "'12.07"
$ws.Range("C11").Value = "'2.57"
$ws.Range("D11").Value = "'14.65"

# Table 1 (Statistical Institution): Enterprises (% of total)
$ws.Range("B13").Value = "'81.83"
$ws.Range("C13").Value = "'17.45"
$ws.Range("D13").Value = "'99.28"

# Table 2 (SME Associations): Enterprises density (per 1000 people)
$ws.Range("D32").Value = "'13.55"
